$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the repeated constant across the full data range B2:K21
$ws.Range("B2:K21").Value = -19.20789370932079

# Override the cells whose value differs from the repeated constant
$ws.Range("C2").Value = 2.435328784732104
$ws.Range("I3").Value = 2.198479999059916
$ws.Range("C4").Value = 2.104127820868193
$ws.Range("D4").Value = 2.875293748672086
$ws.Range("F4").Value = 2.531648359503038
$ws.Range("H4").Value = 1.683015480620963
$ws.Range("J4").Value = 2.119495012731922
$ws.Range("C5").Value = 0.9934752443368481
$ws.Range("G5").Value = 4.321925831565467
$ws.Range("B7").Value = 2.989118057327376
$ws.Range("E8").Value = 2.869639521486216
$ws.Range("B9").Value = 3.592168239083772
$ws.Range("I10").Value = 1.585964366804594
$ws.Range("K10").Value = 1.654269150809091
$ws.Range("E11").Value = 2.003170522279019
$ws.Range("K11").Value = 1.117553395683036
$ws.Range("E13").Value = 1.67682907521837
$ws.Range("J13").Value = 2.522580159446856
$ws.Range("K13").Value = 1.568318415381909
$ws.Range("D14").Value = 1.660643273753569
$ws.Range("K14").Value = 2.26092171805205
$ws.Range("D15").Value = -0.2533083045994791
$ws.Range("J16").Value = 2.405187544033782
$ws.Range("C17").Value = 0.6392814212179742
$ws.Range("D17").Value = -0.07463859005717562
$ws.Range("H17").Value = 0.4320492920493821
$ws.Range("I17").Value = 0.9145058712503737
$ws.Range("J17").Value = 1.392158017242325
$ws.Range("H18").Value = 0.3712767285457253
$ws.Range("I18").Value = 0.919074117223135
$ws.Range("J18").Value = 0.9903846547753683
$ws.Range("D19").Value = 1.722881881187526
$ws.Range("H19").Value = 1.848837266789482
$ws.Range("I19").Value = 2.20092786404715
$ws.Range("C20").Value = 1.703598491457961
$ws.Range("D20").Value = 2.141228208227411
$ws.Range("F20").Value = 3.829611185351949
$ws.Range("H20").Value = 2.414387392804119
$ws.Range("I20").Value = 2.012341124945288
$ws.Range("K20").Value = 2.791602981658981
$ws.Range("C21").Value = 1.801602122106121
$ws.Range("E21").Value = 2.455518895750759
$ws.Range("H21").Value = 2.382202158410758

Write-Output "Applied updated PSSM values"
